$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 172-195 already exist; one weekly record was dropped off the top of the
# 2021-06-to-2023-09 series and a new one appended at the bottom, so every row from
# 172 down now carries the data that used to sit one row below it. Only the
# per-record columns (date, volume, price range, origin, $/kg) actually change --
# the market/region/product columns are identical on every row of this sheet.
$data = @(
    @{Row=172; D=44572; M=170; N=6000; O=6000; P=6000; R="Perú"; S=1500},
    @{Row=173; D=44776; M=228; N=10000; O=10000; P=10000; R="Brasil"; S=2500},
    @{Row=174; D=44315; M=210; N=9000; O=9000; P=9000; R="Perú"; S=2250},
    @{Row=175; D=44432; M=120; N=8000; O=8000; P=8000; R="México"; S=2000},
    @{Row=176; D=44517; M=200; N=6000; O=6000; P=6000; R="Ecuador"; S=1500},
    @{Row=177; D=45027; M=248; N=7000; O=7000; P=7000; R="Perú"; S=1750},
    @{Row=178; D=44781; M=228; N=10000; O=10000; P=10000; R="México"; S=2500},
    @{Row=179; D=44251; M=240; N=5000; O=5000; P=5000; R="Perú"; S=1250},
    @{Row=180; D=44484; M=208; N=7000; O=8000; P=7481; R="Perú"; S=1870},
    @{Row=181; D=44609; M=150; N=7000; O=7000; P=7000; R="Perú"; S=1750},
    @{Row=182; D=44467; M=240; N=8000; O=8000; P=8000; R="Brasil"; S=2000},
    @{Row=183; D=44664; M=210; N=6500; O=6500; P=6500; R="Perú"; S=1625},
    @{Row=184; D=44460; M=200; N=8000; O=8000; P=8000; R="Brasil"; S=2000},
    @{Row=185; D=45012; M=248; N=7000; O=7000; P=7000; R="Perú"; S=1750},
    @{Row=186; D=44222; M=200; N=6000; O=6000; P=6000; R="Perú"; S=1500},
    @{Row=187; D=44362; M=200; N=6000; O=6000; P=6000; R="Perú"; S=1500},
    @{Row=188; D=44826; M=228; N=9000; O=9000; P=9000; R="Brasil"; S=2250},
    @{Row=189; D=44845; M=300; N=7000; O=7000; P=7000; R="Brasil"; S=1750},
    @{Row=190; D=44525; M=120; N=6000; O=6000; P=6000; R="Perú"; S=1500},
    @{Row=191; D=45173; M=250; N=10000; O=10000; P=10000; R="Brasil"; S=2500},
    @{Row=192; D=44200; M=180; N=5000; O=5000; P=5000; R="Perú"; S=1250},
    @{Row=193; D=44991; M=248; N=8000; O=8000; P=8000; R="Perú"; S=2000},
    @{Row=194; D=45191; M=248; N=10000; O=10000; P=10000; R="Brasil"; S=2500},
    @{Row=195; D=44991; M=248; N=8000; O=8000; P=8000; R="Perú"; S=2000}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 13).Value = $item.M
    $ws.Cells.Item($r, 14).Value = $item.N
    $ws.Cells.Item($r, 15).Value = $item.O
    $ws.Cells.Item($r, 16).Value = $item.P
    $ws.Cells.Item($r, 18).Value = $item.R
    $ws.Cells.Item($r, 19).Value = $item.S
}

# Row 196 is brand new (the record that used to live in row 195 before the shift).
# Clone the constant metadata columns from row 195 (same market/product on every
# row of this sheet) using Value2 (plain .Value reads are unreliable here), then
# set the new row's own date/volume/price/origin values.
$newRow = 196
$srcRow = 195
$ws.Cells.Item($newRow, 1).Value = $ws.Cells.Item($srcRow, 1).Value2
$ws.Cells.Item($newRow, 2).Value = $ws.Cells.Item($srcRow, 2).Value2
$ws.Cells.Item($newRow, 3).Value = $ws.Cells.Item($srcRow, 3).Value2
$ws.Cells.Item($newRow, 5).Value = $ws.Cells.Item($srcRow, 5).Value2
$ws.Cells.Item($newRow, 6).Value = $ws.Cells.Item($srcRow, 6).Value2
$ws.Cells.Item($newRow, 7).Value = $ws.Cells.Item($srcRow, 7).Value2
$ws.Cells.Item($newRow, 8).Value = $ws.Cells.Item($srcRow, 8).Value2
$ws.Cells.Item($newRow, 9).Value = $ws.Cells.Item($srcRow, 9).Value2
$ws.Cells.Item($newRow, 10).Value = $ws.Cells.Item($srcRow, 10).Value2
$ws.Cells.Item($newRow, 11).Value = $ws.Cells.Item($srcRow, 11).Value2
$ws.Cells.Item($newRow, 12).Value = $ws.Cells.Item($srcRow, 12).Value2
$ws.Cells.Item($newRow, 17).Value = $ws.Cells.Item($srcRow, 17).Value2
$ws.Cells.Item($newRow, 20).Value = $ws.Cells.Item($srcRow, 20).Value2

$ws.Cells.Item($newRow, 4).Value = 45191
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($srcRow, 4).NumberFormat
$ws.Cells.Item($newRow, 13).Value = 248
$ws.Cells.Item($newRow, 14).Value = 10000
$ws.Cells.Item($newRow, 15).Value = 10000
$ws.Cells.Item($newRow, 16).Value = 10000
$ws.Cells.Item($newRow, 18).Value = "Brasil"
$ws.Cells.Item($newRow, 19).Value = 2500
